# Insert a new weekly record above the current row 75, pushing all
# existing data rows (75-94) down by one (to 76-95), matching the
# "Fruta / hortaliza, semanal" update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at row 75; rows 75:94 shift down to 76:95.
$ws.Rows(75).Insert()

# Populate the new row 75 with this week's record.
$ws.Range("A75").Value = 2
$ws.Range("B75").Value = 'Comercializadora del Agro de Limarí'
$ws.Range("C75").Value = 'Coquimbo'
$ws.Range("D75").Value = 44504
$ws.Range("E75").Value = 4
$ws.Range("F75").Value = 100112024
$ws.Range("G75").Value = 'Choclo'
$ws.Range("H75").Value = 'Dulce o Americano'
$ws.Range("I75").Value = 'Primera'
$ws.Range("J75").Value = 1200
$ws.Range("K75").Value = 37000
$ws.Range("L75").Value = 39000
$ws.Range("M75").Value = 38000
$ws.Range("N75").Value = '$/malla 70 unidades'
$ws.Range("O75").Value = 'Provincia de Limarí'
$ws.Range("P75").Value = 543
$ws.Range("Q75").Value = 70
$ws.Range("R75").Value = 'Hortaliza'
